$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "KL99 removed from database" history row (row 16).
$ws.Range("A16").Value = 45901
$ws.Range("B16").Value = "KL99"
$ws.Range("C16").Value = "removed from database"
$ws.Range("D16").Value = "Frequent mistyping of KL13 as KL99. KL99 differs from KL13 and KL34 by a duplication of  group 154."

# Match the saved view state (scrolled/selected as left by the author).
$ws.Range("C16").Select()
$ws.Application.ActiveWindow.ScrollRow = 12
